$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1145.875
$ws.Range("I2").Value = 655
$ws.Range("K2").Value = 655
$ws.Range("M2").Value = -542
$ws.Range("H10").Value = 40004
$ws.Range("I10").Value = 40004
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 40004
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -39711
$ws.Range("N10").Value = ""
$ws.Range("H15").Value = 774.64703
$ws.Range("I15").Value = 774.64703
$ws.Range("K15").Value = 2323.94109
$ws.Range("M15").Value = -2154.94109
$ws.Range("H17").Value = 1667.6786
$ws.Range("J17").Value = 1678.3334
$ws.Range("L17").Value = 5035.0002
$ws.Range("N17").Value = -5371.0002
$ws.Range("H40").Value = 38468040
$ws.Range("J40").Value = 71434890
$ws.Range("L40").Value = 71434890
$ws.Range("N40").Value = -71435240
$ws.Range("H86").Value = 4433.625
$ws.Range("I86").Value = 2893.8
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 2893.8
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -1770.8
$ws.Range("N86").Value = -9246
$ws.Range("H89").Value = 4433.625
$ws.Range("I89").Value = 2893.8
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 14469
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -8853
$ws.Range("N89").Value = -46232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 200012300
$ws.Range("I36").Value = 333340500
$ws.Range("J36").Value = 19999
$ws.Range("K36").Value = 333340500
$ws.Range("L36").Value = 19999
$ws.Range("M36").Value = -333340154
$ws.Range("N36").Value = -20691
$ws.Range("H102").Value = 8569
$ws.Range("J102").Value = 9636.666999999999
$ws.Range("L102").Value = 9636.666999999999
$ws.Range("N102").Value = -12880.667
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H110").Value = 3909.6553
$ws.Range("J110").Value = 4653.8887
$ws.Range("L110").Value = 4653.8887
$ws.Range("N110").Value = -8743.8887
$ws.Range("H132").Value = 3228084.5
$ws.Range("I132").Value = 2197.0386
$ws.Range("K132").Value = 6591.1158
$ws.Range("M132").Value = -4061.1158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1545
$ws.Range("I22").Value = 1545
$ws.Range("K22").Value = 1545
$ws.Range("M22").Value = -1372
$ws.Range("H105").Value = 696205.6
$ws.Range("I105").Value = 954400.0600000001
$ws.Range("J105").Value = 7687.222
$ws.Range("K105").Value = 954400.0600000001
$ws.Range("L105").Value = 7687.222
$ws.Range("M105").Value = -952653.0600000001
$ws.Range("N105").Value = -11181.222
$ws.Range("H107").Value = 3373.7144
$ws.Range("I107").Value = 3643.3333
$ws.Range("K107").Value = 3643.3333
$ws.Range("M107").Value = -1723.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1470.625
$ws.Range("I16").Value = 541.5
$ws.Range("K16").Value = 541.5
$ws.Range("M16").Value = -254.5
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("H33").Value = 3000
$ws.Range("J33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("N33").Value = -3758
$ws.Range("H86").Value = 9320
$ws.Range("I86").Value = 7080.75
$ws.Range("K86").Value = 7080.75
$ws.Range("M86").Value = -5957.75
$ws.Range("H89").Value = 9320
$ws.Range("I89").Value = 7080.75
$ws.Range("K89").Value = 35403.75
$ws.Range("M89").Value = -29787.75
$ws.Range("H107").Value = 1959.2307
$ws.Range("I107").Value = 495
$ws.Range("K107").Value = 495
$ws.Range("M107").Value = 1425
$ws.Range("H113").Value = 1470.625
$ws.Range("I113").Value = 541.5
$ws.Range("K113").Value = 541.5
$ws.Range("M113").Value = 1628.5
$ws.Range("H132").Value = 2274.7307
$ws.Range("I132").Value = 1912.8948
$ws.Range("K132").Value = 5738.6844
$ws.Range("M132").Value = -3208.6844

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 537.93335
$ws.Range("I5").Value = 301.84616
$ws.Range("J5").Value = 2072.5
$ws.Range("K5").Value = 905.5384799999999
$ws.Range("L5").Value = 6217.5
$ws.Range("M5").Value = -793.5384799999999
$ws.Range("N5").Value = -6441.5
$ws.Range("H62").Value = 12415.9
$ws.Range("I62").Value = 9978
$ws.Range("K62").Value = 29934
$ws.Range("M62").Value = -29248
$ws.Range("H65").Value = 12415.9
$ws.Range("I65").Value = 9978
$ws.Range("K65").Value = 89802
$ws.Range("M65").Value = -86370
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H105").Value = 15871.8
$ws.Range("I105").Value = 10026
$ws.Range("J105").Value = 17333.25
$ws.Range("K105").Value = 30078
$ws.Range("L105").Value = 51999.75
$ws.Range("M105").Value = -27457
$ws.Range("N105").Value = -57241.75
$ws.Range("H121").Value = 8681.111000000001
$ws.Range("J121").Value = 10904.714
$ws.Range("L121").Value = 32714.142
$ws.Range("N121").Value = -35334.142
$ws.Range("H135").Value = 537.93335
$ws.Range("I135").Value = 301.84616
$ws.Range("J135").Value = 2072.5
$ws.Range("K135").Value = 2716.61544
$ws.Range("L135").Value = 18652.5
$ws.Range("M135").Value = -181.61544
$ws.Range("N135").Value = -23722.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10674.625
$ws.Range("I70").Value = 8580.799999999999
$ws.Range("J70").Value = 12522.117
$ws.Range("K70").Value = 8580.799999999999
$ws.Range("L70").Value = 12522.117
$ws.Range("M70").Value = -8310.799999999999
$ws.Range("N70").Value = -13062.117
$ws.Range("H73").Value = 10674.625
$ws.Range("I73").Value = 8580.799999999999
$ws.Range("J73").Value = 12522.117
$ws.Range("K73").Value = 8580.799999999999
$ws.Range("L73").Value = 12522.117
$ws.Range("M73").Value = -7644.799999999999
$ws.Range("N73").Value = -14394.117
$ws.Range("H107").Value = 894.25
$ws.Range("I107").Value = 192.33333
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 192.33333
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1727.66667
$ws.Range("N107").Value = -6840
$ws.Range("H126").Value = 4192.222
$ws.Range("I126").Value = 4130.7856
$ws.Range("J126").Value = 4407.25
$ws.Range("K126").Value = 12392.3568
$ws.Range("L126").Value = 13221.75
$ws.Range("M126").Value = -9922.356800000001
$ws.Range("N126").Value = -18161.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6682.769
$ws.Range("I22").Value = 12089.4
$ws.Range("J22").Value = 3303.625
$ws.Range("K22").Value = 12089.4
$ws.Range("L22").Value = 3303.625
$ws.Range("M22").Value = -11794.4
$ws.Range("N22").Value = -3893.625
$ws.Range("H27").Value = 6682.769
$ws.Range("I27").Value = 12089.4
$ws.Range("J27").Value = 3303.625
$ws.Range("K27").Value = 12089.4
$ws.Range("L27").Value = 3303.625
$ws.Range("M27").Value = -11982.4
$ws.Range("N27").Value = -3517.625
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = ""
$ws.Range("H132").Value = 3354.8333
$ws.Range("I132").Value = 2144.7058
$ws.Range("J132").Value = 5412.05
$ws.Range("K132").Value = 6434.117400000001
$ws.Range("L132").Value = 16236.15
$ws.Range("M132").Value = -3904.117400000001
$ws.Range("N132").Value = -21296.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5533
$ws.Range("I107").Value = 3916.5
$ws.Range("K107").Value = 11749.5
$ws.Range("M107").Value = -9829.5
$ws.Range("H109").Value = 120449.5
$ws.Range("J109").Value = 120449.5
$ws.Range("L109").Value = 120449.5
$ws.Range("N109").Value = -123223.5
$ws.Range("H126").Value = 4393.4116
$ws.Range("I126").Value = 5509
$ws.Range("J126").Value = 2799.7144
$ws.Range("K126").Value = 16527
$ws.Range("L126").Value = 8399.143199999999
$ws.Range("M126").Value = -14057
$ws.Range("N126").Value = -13339.1432
